$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 4 new rows before row 12 (shifts old rows 12-40 down to 16-44).
#    Copy formatting (incl. borders) from row 11 onto the freshly inserted
#    rows so they keep the same style indices as the rest of the table
#    instead of the engine's default "no border" insert formatting.
# ---------------------------------------------------------------------------
$ws.Rows("12:15").Insert() | Out-Null
$ws.Range("A11:F11").Copy() | Out-Null
$ws.Range("A12:F15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Row heights / font-scale (dyDescent is derived from font size, which we
#    don't touch) for the rows whose height changed.
# ---------------------------------------------------------------------------
$ws.Rows("5").RowHeight = 93
$ws.Rows("6").RowHeight = 49.2
$ws.Rows("12").RowHeight = 61.2
$ws.Rows("13").RowHeight = 25.2
$ws.Rows("14").RowHeight = 17.4
$ws.Rows("15").RowHeight = 17.4

# ---------------------------------------------------------------------------
# 3) New content: "Feature 3: Sizing limits" entry (row 12).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Feature 3: Sizing limits"
$ws.Range("B12").Value = 0.5
$ws.Range("C12").Value = "18/05/2025"
$ws.Range("D12").Value = "I have made several updates to the mazeSettings in order to prevent accidentally assigning excessively large values from Unity. Additionally, I have implemented additional checks in the setters to ensure proper validation. The view now utilizes a slider, and upon starting the application, I automatically set the minimum and maximum values to align with the sizing limits specified in MazeSettings"

# ---------------------------------------------------------------------------
# 4) New content: "Feature  3: Improving the UX" entry (row 13).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Feature  3: Improving the UX"
$ws.Range("B13").Value = 0.5
$ws.Range("C13").Value = "18/05/2025"
$ws.Range("D13").Value = "Current slider does not show actual values, therefore user lacks controls. I decided to add 2 an extra text input where user can specify the  size more exactly."
$ws.Range("E13").Value = "X"

# Rows 14 and 15 stay blank placeholder rows (same as the other blank rows in
# the table) - nothing else to set there.

# ---------------------------------------------------------------------------
# 5) Update the "Total amount of hours" formula (now on row 41) to include
#    the 4 extra rows (old range E4:E35 / B4:B35 -> E4:E39 / B4:B39).
# ---------------------------------------------------------------------------
$ws.Range("B41").Formula = '=SUMIF(E4:E39,"<>x",B4:B39)'

# ---------------------------------------------------------------------------
# 6) Sheet view bookkeeping (scroll position / selection) to match the saved
#    state in the workbook.
# ---------------------------------------------------------------------------
$ws.Range("D15").Select()
$excel.ActiveWindow.ScrollRow = 8
